$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 3 content
$ws.Range("B3").Value = "FirefoxProfile"

$code = "about:config`nFirefoxProfile profile = new FireFoxprofile();`nprofile.setPreference(`"`",`"`");`ncapability.setCapability(FireFoxDriver.Profile,profile);"
$ws.Range("C3").Value = $code

# Formatting for C3: wrap text + specific font/color
$ws.Range("C3").WrapText = $true
$ws.Range("C3").Font.Name = "微软雅黑"
$ws.Range("C3").Font.Size = 11
$ws.Range("C3").Font.Color = 2236962

# Row height for row 3
$ws.Rows.Item(3).RowHeight = 66

# Selection change
$ws.Range("D3").Select()
